$d = $word.ActiveDocument

$d.Content.Find.Execute("18×13=234", $true, $false, $false, $false, $false, $true, 1, $false, "99×38=3762", 2) | Out-Null
$d.Content.Find.Execute("41×43=1763", $true, $false, $false, $false, $false, $true, 1, $false, "81×70=5670", 2) | Out-Null
$d.Content.Find.Execute("45×58=2610", $true, $false, $false, $false, $false, $true, 1, $false, "38×37=1406", 2) | Out-Null
$d.Content.Find.Execute("87×46=4002", $true, $false, $false, $false, $false, $true, 1, $false, "79×19=1501", 2) | Out-Null
$d.Content.Find.Execute("52×68=3536", $true, $false, $false, $false, $false, $true, 1, $false, "14×75=1050", 2) | Out-Null
$d.Content.Find.Execute("27×67=1809", $true, $false, $false, $false, $false, $true, 1, $false, "58×65=3770", 2) | Out-Null
$d.Content.Find.Execute("58×95=5510", $true, $false, $false, $false, $false, $true, 1, $false, "35×33=1155", 2) | Out-Null
$d.Content.Find.Execute("28×88=2464", $true, $false, $false, $false, $false, $true, 1, $false, "77×73=5621", 2) | Out-Null
$d.Content.Find.Execute("84×37=3108", $true, $false, $false, $false, $false, $true, 1, $false, "69×83=5727", 2) | Out-Null
$d.Content.Find.Execute("60×98=5880", $true, $false, $false, $false, $false, $true, 1, $false, "34×94=3196", 2) | Out-Null
$d.Content.Find.Execute("46×15=690", $true, $false, $false, $false, $false, $true, 1, $false, "71×61=4331", 2) | Out-Null
$d.Content.Find.Execute("90×37=3330", $true, $false, $false, $false, $false, $true, 1, $false, "33×50=1650", 2) | Out-Null
$d.Content.Find.Execute("27×18=486", $true, $false, $false, $false, $false, $true, 1, $false, "78×88=6864", 2) | Out-Null
$d.Content.Find.Execute("43×79=3397", $true, $false, $false, $false, $false, $true, 1, $false, "64×22=1408", 2) | Out-Null
$d.Content.Find.Execute("61×58=3538", $true, $false, $false, $false, $false, $true, 1, $false, "43×52=2236", 2) | Out-Null
$d.Content.Find.Execute("34×46=1564", $true, $false, $false, $false, $false, $true, 1, $false, "72×44=3168", 2) | Out-Null
$d.Content.Find.Execute("13×19=247", $true, $false, $false, $false, $false, $true, 1, $false, "67×56=3752", 2) | Out-Null
$d.Content.Find.Execute("77×80=6160", $true, $false, $false, $false, $false, $true, 1, $false, "97×71=6887", 2) | Out-Null
$d.Content.Find.Execute("41×37=1517", $true, $false, $false, $false, $false, $true, 1, $false, "24×34=816", 2) | Out-Null
$d.Content.Find.Execute("22×26=572", $true, $false, $false, $false, $false, $true, 1, $false, "17×42=714", 2) | Out-Null
$d.Content.Find.Execute("18×86=1548", $true, $false, $false, $false, $false, $true, 1, $false, "79×52=4108", 2) | Out-Null
$d.Content.Find.Execute("42×82=3444", $true, $false, $false, $false, $false, $true, 1, $false, "85×65=5525", 2) | Out-Null
$d.Content.Find.Execute("93×28=2604", $true, $false, $false, $false, $false, $true, 1, $false, "22×61=1342", 2) | Out-Null
$d.Content.Find.Execute("57×50=2850", $true, $false, $false, $false, $false, $true, 1, $false, "95×88=8360", 2) | Out-Null
$d.Content.Find.Execute("98×40=3920", $true, $false, $false, $false, $false, $true, 1, $false, "82×60=4920", 2) | Out-Null
